$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FBS")
$ws.Range("AK2").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK3").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK4").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK5").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK6").Value = "2024-10-11T10:01:39.733003"
$ws.Range("M6").Value = "ESE"
$ws.Range("O6").Value = 85.79000000000001
$ws.Range("P6").Value = 5.8
$ws.Range("Q6").Value = "ESE"
$ws.Range("S6").Value = -0.72
$ws.Range("T6").Value = -0.72
$ws.Range("U6").Value = -0.9
$ws.Range("AE7").Value = 0.03149606299212598
$ws.Range("AK7").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Y7").Value = 65.5
$ws.Range("Z7").Value = -115
$ws.Range("A8").Value = "Boise State @ Hawaii"
$ws.Range("AA8").Value = $null
$ws.Range("AB8").Value = $null
$ws.Range("AE8").Value = -0.01626016260162602
$ws.Range("AF8").Value = $null
$ws.Range("AK8").Value = "2024-10-11T10:01:39.733003"
$ws.Range("C8").Value = "05:00 PM"
$ws.Range("E8").Value = "NW-SE"
$ws.Range("F8").Value = "High"
$ws.Range("G8").Value = $null
$ws.Range("H8").Value = -803.2221069
$ws.Range("I8").Value = 75.33
$ws.Range("J8").Value = 53.65
$ws.Range("K8").Value = 13.5
$ws.Range("L8").Value = 2014
$ws.Range("M8").Value = "W"
$ws.Range("N8").Value = "WSW"
$ws.Range("O8").Value = 83.72
$ws.Range("P8").Value = 5.4
$ws.Range("Q8").Value = "W"
$ws.Range("S8").Value = -0.46
$ws.Range("T8").Value = -0.46
$ws.Range("U8").Value = -8.1
$ws.Range("V8").Value = "21.294294, -157.819338"
$ws.Range("W8").Value = 61.5
$ws.Range("Y8").Value = 60.5
$ws.Range("Z8").Value = -110
$ws.Range("A9").Value = "San Jose State @ Colorado State"
$ws.Range("AA9").Value = 2
$ws.Range("AB9").Value = -1.5
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 3.5
$ws.Range("AK9").Value = "2024-10-11T10:01:39.733003"
$ws.Range("C9").Value = "01:30 PM"
$ws.Range("E9").Value = "N-S"
$ws.Range("F9").Value = "Med"
$ws.Range("G9").Value = "E/W"
$ws.Range("H9").Value = 1500.09291306
$ws.Range("I9").Value = 48.81
$ws.Range("J9").Value = 59.64
$ws.Range("K9").Value = 8.1
$ws.Range("L9").Value = 2017
$ws.Range("M9").Value = "SW"
$ws.Range("N9").Value = "NNW"
$ws.Range("O9").Value = 82.01000000000001
$ws.Range("P9").Value = 6.6
$ws.Range("Q9").Value = "NNW"
$ws.Range("S9").Value = -0.25
$ws.Range("T9").Value = -3.5
$ws.Range("U9").Value = -1.5
$ws.Range("V9").Value = "40.570015, -105.088435"
$ws.Range("W9").Value = 55.5
$ws.Range("Y9").Value = 55.5
$ws.Range("Z9").Value = -105
$ws.Range("AK10").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK11").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q11").Value = "ESE"
$ws.Range("Z11").Value = -105
$ws.Range("AB12").Value = -21.5
$ws.Range("AF12").Value = 3
$ws.Range("AK12").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK13").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK14").Value = "2024-10-11T10:01:39.733003"
$ws.Range("O14").Value = 55.52
$ws.Range("AK15").Value = "2024-10-11T10:01:39.733003"
$ws.Range("M15").Value = "W"
$ws.Range("N15").Value = "WNW"
$ws.Range("O15").Value = 69.32000000000001
$ws.Range("P15").Value = 4.7
$ws.Range("Q15").Value = "W"
$ws.Range("U15").Value = -0.8
$ws.Range("AK16").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q16").Value = "S"
$ws.Range("AK17").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK18").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK19").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK20").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK21").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK22").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK23").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK24").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK25").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK26").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK27").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK28").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK29").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK30").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK31").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK32").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK33").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK34").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK35").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK36").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK37").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK38").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK39").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK40").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK41").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK42").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q42").Value = "W"
$ws.Range("AK43").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK44").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK45").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK46").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q46").Value = "NW"
$ws.Range("AK47").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK48").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q48").Value = "S"
$ws.Range("AK49").Value = "2024-10-11T10:01:39.733003"
$ws.Range("AK50").Value = "2024-10-11T10:01:39.733003"
$ws.Range("Q50").Value = "N"

$ws = $wb.Worksheets.Item("Other")
$ws.Range("O19").Value = "ENE"
$ws.Range("P19").Value = "NE"
$ws.Range("Q19").Value = 58.58000000000001
$ws.Range("R19").Value = 5.9
$ws.Range("S42").Value = "ESE"
